{"js": "// Add two new \"Algorithms and data structures\" note paragraphs right after the\n// \"Keyboard shortcuts\" hyperlink paragraph, an empty spacer paragraph, and a\n// new leading sentence on the final (bookmarked) paragraph.\n//\n// We build the new content as a Flat-OPC OOXML fragment and insert it with\n// Range.insertOoxml() so the exact run layout (including the bold \"S \" run\n// and the \"an\"-wrapping grammar-check markers the author's Word session\n// produced) is reproduced faithfully instead of being collapsed/merged the\n// way plain insertText() calls would be.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The very last paragraph in the body is the one holding the `_GoBack`\n// bookmark \u2014 the new content is spliced in immediately before it (the last\n// new run lands inside that same paragraph, ahead of the bookmark).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst flatOpcPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">To scale </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  \"<w:r><w:t>an</w:t></w:r>\" +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> mesh, press </w:t></w:r>' +\n  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">S </w:t></w:r>' +\n  \"<w:r><w:t>and the object will start to scale in all directions based on your mouse movements.</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>Control c</w:t></w:r>\" +\n  \"<w:r><w:t>hanges it to</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> incremental steps</w:t></w:r>' +\n  \"<w:r><w:t>, while shift makes subtle movements.</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">After adding a new object, for some objects, you will see an options bar to the side that allows you to determine the </w:t></w:r>' +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nlastParagraph.getRange(\"Start\").insertOoxml(flatOpcPackage, \"Before\");\nawait context.sync();\n", "ps1": "# Add two new \"Algorithms and data structures\" note paragraphs right after\n# the \"Keyboard shortcuts\" hyperlink paragraph, an empty spacer paragraph,\n# and a new leading sentence on the final (bookmarked) paragraph.\n#\n# The new content is built as a Flat-OPC OOXML fragment and spliced in with\n# Range.InsertXML() at a zero-length Range collapsed to the very start of\n# the document's last paragraph. InsertXML() *replaces* the content of the\n# range it is called on, and a zero-length range sits exactly at that\n# insertion point, so it behaves like an \"insert before\" \u2014 it reproduces the\n# exact run layout (the bold \"S \" run and the \"an\"-wrapping grammar-check\n# markers the author's Word session produced) instead of merging adjacent\n# same-format runs the way Range.InsertAfter(text) would, and it leaves the\n# `_GoBack` bookmark and the paragraph's own identity untouched.\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$insertionPoint = $d.Range($lastParagraph.Range.Start, $lastParagraph.Range.Start)\n\n$flatOpcPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">To scale </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>an</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> mesh, press </w:t></w:r>' +\n    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">S </w:t></w:r>' +\n    '<w:r><w:t>and the object will start to scale in all directions based on your mouse movements.</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p>' +\n    '<w:r><w:t>Control c</w:t></w:r>' +\n    '<w:r><w:t>hanges it to</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> incremental steps</w:t></w:r>' +\n    '<w:r><w:t>, while shift makes subtle movements.</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p/>' +\n    '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">After adding a new object, for some objects, you will see an options bar to the side that allows you to determine the </w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$insertionPoint.InsertXML($flatOpcPackage)\n"}
